$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    ", and C# .NET, familiar",  # FindText
    $false,                      # MatchCase
    $false,                      # MatchWholeWord
    $false,                      # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    " and familiar",             # ReplaceWith
    2                            # Replace (wdReplaceAll)
)
